$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "accuracy_balanced" metric (mean + std) based on reviewer
# feedback. This inserts a new 7-row "accuracy_balanced_mean" block right
# after the existing "f1_micro_mean" block (pushing the two *_std blocks
# down by 7 rows), and appends a new 7-row "accuracy_balanced_std" block
# at the end of the table.
# ---------------------------------------------------------------------------

# 1) Insert 7 blank rows for the new accuracy_balanced_mean block, right
#    before the current f1_macro_std block (old row 16).
$ws.Range("A16:A22").EntireRow.Insert()

# Copy the label-column formatting (bold, centered, bordered) from an
# existing data row onto the freshly inserted rows (column A only).
$ws.Range("A2").Copy()
$ws.Range("A16:A22").PasteSpecial(-4122)

# 2) Fill in the new accuracy_balanced_mean rows (now at 16-22).
$meanRows = @(
    @("0",     0,     0,     0,     0,     0,     0.446),
    @("100",   0.173, 0.12,  0.173, 0.185, 0.216, 0.525),
    @("500",   0.302, 0.266, 0.297, 0.332, 0.472, 0.598),
    @("1000",  0.346, 0.348, 0.348, 0.401, 0.524, 0.644),
    @("2500",  0.375, 0.406, 0.395, 0.489, 0.626, 0.668),
    @("3000",  0.378, 0.427, 0.429, 0.54,  0.681, 0.707),
    @("10000", 0.428, 0.402, 0.447, 0.572, 0.691, 0.722)
)

$r = 16
foreach ($row in $meanRows) {
    $ws.Cells.Item($r, 1).Value = "accuracy_balanced_mean"
    $ws.Cells.Item($r, 2).Value = "'" + $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# 3) Append the new accuracy_balanced_std block at the end of the table
#    (rows 37-43).
$stdRows = @(
    @("0",     0,     0,     0,     0,     0,     0),
    @("100",   0.001, 0.007, 0.007, 0.006, 0.023, 0.014),
    @("500",   0.02,  0.014, 0.009, 0.017, 0.02,  0.014),
    @("1000",  0.009, 0.004, 0.01,  0.016, 0.018, 0.004),
    @("2500",  0.001, 0.001, 0.01,  0.007, 0.024, 0.014),
    @("3000",  0.003, 0.004, 0.008, 0.005, 0.006, 0.003),
    @("10000", 0.004, 0.003, 0.005, 0.009, 0.02,  0.013)
)

# Apply the label-column formatting to the new rows first.
$ws.Range("A2").Copy()
$ws.Range("A37:A43").PasteSpecial(-4122)

$r = 37
foreach ($row in $stdRows) {
    $ws.Cells.Item($r, 1).Value = "accuracy_balanced_std"
    $ws.Cells.Item($r, 2).Value = "'" + $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $r++
}
